# Updated symbol list on Sat Dec 31 04:39:39 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    # Force text storage (avoids Excel auto-converting numeric-looking
    # strings like "245.14" into numbers), then reset the style back to
    # Normal so no extra cell style/number-format gets persisted.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Column D price updates (rows unaffected by reordering)
Set-TextCell "D2"  "245.14"
Set-TextCell "D3"  "25.35"
Set-TextCell "D4"  "5.120"
Set-TextCell "D5"  "0.05573"
Set-TextCell "D6"  "6.499"
Set-TextCell "D7"  "3.019"
Set-TextCell "D8"  "0.8173"
Set-TextCell "D9"  "0.8465"
Set-TextCell "D10" "0.1345"
Set-TextCell "D11" "0.06958"

# Rows 12-20: coin list re-ranked, shifting entries up and wrapping row 12's
# original coin (Liechtenstein Cryptoassets Exchange) down to row 20.
Set-TextCell "B12" "BitrueCoin"
Set-TextCell "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.02885"
Set-TextCell "E12" "11BitrueCoinBTR"

Set-TextCell "B13" "BitMartToken"
Set-TextCell "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D13" "0.09377"
Set-TextCell "E13" "12BitMartTokenBMX"

Set-TextCell "B14" "BitForexToken"
Set-TextCell "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D14" "0.001514"
Set-TextCell "E14" "13BitForexTokenBF"

Set-TextCell "B15" "TigerCash"
Set-TextCell "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D15" "0.006131"
Set-TextCell "E15" "14TigerCashTCH"

Set-TextCell "B16" "LEO"
Set-TextCell "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D16" "3.500"
Set-TextCell "E16" "15LEOLEO"

Set-TextCell "B17" "BTSEToken"
Set-TextCell "C17" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D17" "2.063"
Set-TextCell "E17" "16BTSETokenBTSE"

Set-TextCell "B18" "One"
Set-TextCell "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D18" "0.009823"
Set-TextCell "E18" "17OneONE"

Set-TextCell "B19" "BitpandaEcosystemToken"
Set-TextCell "C19" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell "D19" "0.3179"
Set-TextCell "E19" "18BitpandaEcosystemTokenBEST"

Set-TextCell "B20" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C20" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D20" "0.03159"
Set-TextCell "E20" "19LiechtensteinCryptoassetsExchangeLCX"

# Remaining scattered price / label updates
Set-TextCell "D22" "3.742"
Set-TextCell "D23" "0.04709"
Set-TextCell "D24" "0.1375"
Set-TextCell "D25" "0.001247"

Set-TextCell "D27" "0.00009703"
Set-TextCell "E27" "26NitroExNTXBestin24h"

Set-TextCell "D41" "0.006212"

Set-TextCell "D43" "0.002383"
Set-TextCell "D44" "0.008328"
Set-TextCell "D45" "0.00005295"

Set-TextCell "D47" "0.1500"
Set-TextCell "D48" "0.002122"
